# Generate Report for Handback
# Update the "Correspond Handoff/Handback Datetime" and "Latest HO Xliff
# Generate Date" timestamps for the 8fdaf50e-...md file row (row 3) on the
# Overview, zh-cn and de-de sheets, reflecting a freshly regenerated
# handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: G3 = "Latest HO Xliff Generate Date" for 8fdaf50e-....md
$overview.Range("G3").Value = "2016-08-31 12:55:54"

# zh-cn sheet: row 3 corresponds to 8fdaf50e-....md
$zhcn.Range("H3").Value = "2016-08-31 12:55:49"   # Correspond Handoff Datetime
$zhcn.Range("K3").Value = "2016-08-31 12:56:32"   # Correspond Handback DateTime

# de-de sheet: row 3 corresponds to 8fdaf50e-....md
$dede.Range("H3").Value = "2016-08-31 12:55:54"   # Correspond Handoff Datetime (matches Overview G3)
$dede.Range("K3").Value = "2016-08-31 12:56:39"   # Correspond Handback DateTime
